$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.078.66"
$ws.Range("E2").Value = "  +4.54%  "
$ws.Range("D3").Value = "3.465.29"
$ws.Range("E3").Value = "  +4.18%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'585.48"
$ws.Range("E5").Value = "  +6.17%  "
$ws.Range("D6").Value = "'187.81"
$ws.Range("E6").Value = "  +8.61%  "
$ws.Range("D7").Value = "'0.633"
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("D8").Value = "3.459.41"
$ws.Range("E8").Value = "  +4.30%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").Value = "'0.647"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").Value = "'56.67"
$ws.Range("E12").Value = "  +6.41%  "
$ws.Range("D13").Value = "'0.0000278"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "'9.42"
$ws.Range("E14").Value = "  +4.02%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'18.74"
$ws.Range("E15").Value = "  +3.67%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.472.59"
$ws.Range("E16").Value = "  +4.75%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "67.140.41"
$ws.Range("E17").Value = "  +4.71%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'12.16"
$ws.Range("E18").Value = "  +3.97%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.118"
$ws.Range("E19").Value = "  -1.72%  "
$ws.Range("B20").Value = "Polygon"
$ws.Range("C20").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D20").Value = "'1.02"
$ws.Range("E20").Value = "  +3.72%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'486.09"
$ws.Range("E21").Value = "  +8.81%  "
$ws.Range("B22").Value = "Toncoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D22").Value = "'5.38"
$ws.Range("E22").Value = "  +8.50%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "'16.83"
$ws.Range("E23").Value = "  +22.13%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'4.46"
$ws.Range("E24").Value = "  +10.65%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'89.66"
$ws.Range("E25").Value = "  +3.28%  "
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").Value = "'2.94"
$ws.Range("E26").Value = "  +2.80%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'10.93"
$ws.Range("E27").Value = "  +3.17%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "'9.08"
$ws.Range("E28").Value = "  +6.19%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'31.38"
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'7.14"
$ws.Range("E30").Value = "  +10.01%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'601.65"
$ws.Range("E31").Value = "  +5.50%  "
$ws.Range("D32").Value = "'11.74"
$ws.Range("E32").Value = "  +3.42%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").Value = "'64.27"
$ws.Range("E33").Value = "  +3.34%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.111"
$ws.Range("E34").Value = "  +5.05%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.150"
$ws.Range("E35").Value = "  +6.67%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "'36.58"
$ws.Range("E37").Value = "  +4.08%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'3.52"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "'0.384"
$ws.Range("E39").Value = "  +4.97%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0756"
$ws.Range("E40").Value = "  +3.88%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.234.09"
$ws.Range("E41").Value = "  +5.75%  "
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").Value = "'2.90"
$ws.Range("E42").Value = "  +6.77%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0429"
$ws.Range("E43").Value = "  +4.25%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.80"
$ws.Range("E44").Value = "  +24.91%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.27"
$ws.Range("E45").Value = "  +3.07%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.51"
$ws.Range("E46").Value = "  +3.21%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.135"
$ws.Range("E47").Value = "  +1.55%  "
$ws.Range("B48").Value = "LidoDAOToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D48").Value = "'3.27"
$ws.Range("E48").Value = "  +12.68%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").Value = "'8.69"
$ws.Range("E50").Value = "  +6.63%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'139.73"
$ws.Range("E51").Value = "  -1.81%  "
